$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in column-name cell and correct the data type for the ringi status row
$ws.Range("A7").Value = "ringiStatus"
$ws.Range("C7").Value = "varchar(255)"

# Remove the stray trailing space value from G9
$ws.Range("G9").ClearContents()

# Remove the stray leftover row 10 entirely
$ws.Rows("10:10").Delete()

# Move the active selection to C8
$ws.Range("C8").Select()
